# Update the simulation results table: keep only the two L-Glutamate
# exchange rows (re-computed values), dropping the Nitrate/Glycine/Citrate
# rows entirely so the sheet shrinks from A1:E8 to A1:E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 4-8 (Glycine__Ibark, Glycine__Phellogen, L-Glutamate__Leaf,
# L-Glutamate__Phellogen, Citrate__Leaf) - Delete() shifts the used range
# up so the sheet dimension shrinks accordingly.
$ws.Range("A4:E8").Delete()

# Row 2: was Nitrate__Ibark_Day_sp_exchange -> becomes the (re-run)
# L-Glutamate__Leaf_Day_sp_exchange values.
$ws.Range("A2").Value = "L-Glutamate__Leaf_Day_sp_exchange"
$ws.Range("B2").Value = -0.15
$ws.Range("C2").Value = -44.84272613958672
$ws.Range("D2").Value = 43.38661379847885
$ws.Range("E2").Value = $false

# Row 3: was Glycine__Leaf_Day_sp_exchange -> becomes the (re-run)
# L-Glutamate__Phellogen_Day_sp_exchange values.
$ws.Range("A3").Value = "L-Glutamate__Phellogen_Day_sp_exchange"
$ws.Range("B3").Value = -0.0277777777777779
$ws.Range("C3").Value = -41.17279541871761
$ws.Range("D3").Value = 45.2925391129859
$ws.Range("E3").Value = $false
